$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.6262526512146
$ws.Range("B1").Value = 2.653408527374268
$ws.Range("C1").Value = 5.493977069854736
$ws.Range("D1").Value = 2.810638427734375
$ws.Range("E1").Value = 0.8683362603187561
